# Append two new observation rows (rows 3 and 4) to the BIIB Noun sheet,
# matching the rows already present for the "Noun" method, and widen
# column A slightly to accommodate the new date/time values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 42600.835023148145
$ws.Range("B3").Value = "Noun"
$ws.Range("C3").Value = 7281
$ws.Range("D3").Value = 6046
$ws.Range("E3").Value = 1110
$ws.Range("F3").Value = 147
$ws.Range("G3").Value = 52
$ws.Range("H3").Value = 73
$ws.Range("I3").Value = 26
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 99
$ws.Range("M3").Value = 0

# Row 4
$ws.Range("A4").Value = 42600.879988425928
$ws.Range("B4").Value = "Noun"
$ws.Range("C4").Value = 7788
$ws.Range("D4").Value = 6213
$ws.Range("E4").Value = 1124
$ws.Range("F4").Value = 152
$ws.Range("G4").Value = 57
$ws.Range("H4").Value = 71
$ws.Range("I4").Value = 26
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 99
$ws.Range("M4").Value = 0

# Column A (the date column) widens slightly now that the new timestamps
# render one character wider than the existing one.
$ws.Columns.Item(1).ColumnWidth = 14
